$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that sits on the paragraph
#    ending in "(先不做)" (4th paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Merge the "输出：导通损耗、开关损耗、效率" paragraph (2nd occurrence)
#    with the empty paragraph that follows it, by deleting the paragraph
#    mark between them.
$outputPara = $d.Paragraphs.Item(6)
$mergeStart = $outputPara.Range.End - 1
$mergeRange = $d.Range($mergeStart, $outputPara.Range.End)
$mergeRange.Delete()

# 3. Give that (now merged) paragraph's mark an eastAsia font hint, same
#    as the rest of the document's paragraph marks.
$outputPara = $d.Paragraphs.Item(6)
$outputPara.Range.Font.NameFarEast = "eastAsia"

# 4. Re-create the "_GoBack" bookmark at the end of that paragraph (after
#    its run, before the paragraph mark).
$bmPos = $outputPara.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
